$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Journal de travail")

# Fix the previously blank description on row 65 (was a placeholder single space)
$ws.Range("D65").Value = "Refactor"

# Add the new journal entry on row 66
$ws.Range("A66").Value = 45110
$ws.Range("B66").Value = "Implémentation"
$ws.Range("C66").Value = 4
$ws.Range("D66").Value = "Réécriture des flux live des écrans avec websockets"

# Move the active selection to D67 (as recorded in the saved view state)
$ws.Range("D67").Select()
